# Insert "TypeScript," as three new runs (" ", "TypeScript", ",") between
# the existing runs " React Hooks," and " JavaScript" in the "Front End:"
# technical-skills line, so it reads:
#   React, React Hooks, TypeScript, JavaScript (ES5/ES6), ...
#
# This runtime rebuilds/coalesces a paragraph's *trailing* runs whenever an
# edit (Find-replace, InsertBefore/After, Range.Text assignment, ...) lands
# strictly before the paragraph's final run, which would silently destroy
# the run layout (and the w:rsidR/... attributes) of every run after the
# insertion point. Inserting exactly at the paragraph's current end avoids
# that. So: Cut the existing tail (everything from right after
# " React Hooks," through the end of the line) to the clipboard, type the
# three new runs at the (now) paragraph end, and Paste the tail back - the
# clipboard Paste reproduces the original run XML byte-for-byte (including
# rsid attributes), unlike re-typing the text.

$d = $word.ActiveDocument

function Find-Range($text) {
    $r = $d.Content
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $text"
    }
    return $r
}

# " React Hooks," is unique to the "Front End:" skills line (the document's
# other "React Hooks" mentions read differently around them), and
# "Materialize" is the last word of that same line / unique in the doc, so
# together they unambiguously bound the tail to relocate.
$hooksEnd = Find-Range("React Hooks,").End
$lineEnd = Find-Range("Materialize").End

$tail = $d.Range($hooksEnd, $lineEnd)
$tail.Cut()

# Type the three new runs at the paragraph's new (current) end, formatting
# each to match the surrounding Calibri runs.
$anchor = "React Hooks,"
foreach ($piece in @(" ", "TypeScript", ",")) {
    $ins = Find-Range($anchor)
    $ins.Collapse(0)
    $ins.InsertAfter($piece)

    $fmt = Find-Range($anchor)
    $fmt.Collapse(0)
    $fmt.MoveEnd(1, $piece.Length)
    $fmt.Font.Name = "Calibri"
    $fmt.Font.NameFarEast = "Calibri"
    $fmt.Font.NameBi = "Calibri"

    $anchor = $anchor + $piece
}

# Paste the original tail back at the paragraph's new end.
$dest = Find-Range($anchor)
$dest.Collapse(0)
$dest.Paste()
